$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set every Fitness value (column C, rows 2-252) to the constant 7293
$ws.Range("C2:C252").Value = 7293
